$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6366130709648132
$ws.Range("B1").Value = 4.171104431152344
$ws.Range("C1").Value = 4.070846557617188
$ws.Range("D1").Value = 1.506852626800537
$ws.Range("E1").Value = 1.079263091087341
